# Inbox data laboproef aangevuld op 04/06/18
# Append a new logger-reading row (date 04/06/2018) to the "Blad1" overview
# sheet. Row 4's B:F cells already carry the correct table formatting; only
# the date cell (A4) needs its number format promoted to match A3 (the
# previous data row) before the value is written. Finish with the
# selection on F5 (the next empty row), matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Activate()

# Copy A3's format (date number format + border) down onto A4.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = 43255
$ws.Range("B4").Value = 4435
$ws.Range("C4").Value = 5370
$ws.Range("D4").Value = 5369
$ws.Range("E4").Value = 4435
$ws.Range("F4").Value = 4435

$ws.Range("F5").Select()
